# Finished Week 13 logging
$wb = $excel.ActiveWorkbook

# OFF sheet: row 3 (R) - Short Att and Short Comp increased
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 119
$wsOff.Range("C3").Value = 88

# DEF sheet: row 3 (R) - Short Att, Short Comp, Deep Att, Deep Comp all increased
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 210
$wsDef.Range("C3").Value = 135
$wsDef.Range("D3").Value = 59
$wsDef.Range("E3").Value = 21
